# Registrar avance del día: se actualizan las horas consumidas (Día 1, columna H)
# de las tareas de la hoja "Casos de Uso".
#  - Fila 6  (tarea de Víctor): se retira el registro de 1 hora consumida.
#  - Fila 7  (tarea "Modulo de pagos ... pago de alumno"): se registran 2 horas
#    consumidas al implementar la primera parte del caso de uso de pago de alumno.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Quitar el valor previamente registrado en H6 (vuelve a quedar vacío)
$ws.Range("H6").ClearContents()

# Registrar 2 horas consumidas en H7
$ws.Range("H7").Value = 2
